$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16169548034668
$ws.Range("B1").Value = 2.414951801300049
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.380393505096436
$ws.Range("E1").Value = 1.231752872467041
